# Natmi following Dr Hou advice
# Update Fgf7-Fgfr1 LR-pair table: sending/target clusters now include ECs, FAPs, sCs
# (previously only FAPs and sCs were sending clusters); recompute downstream metrics.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Fgf7"
$ws.Cells.Item(2,3).Value = "Fgfr1"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 2
$ws.Cells.Item(2,6).Value = 0.6666666666666666
$ws.Cells.Item(2,7).Value = 0.245245
$ws.Cells.Item(2,8).Value = 0.735735
$ws.Cells.Item(2,9).Value = 0.0130094690177091
$ws.Cells.Item(2,10).Value = 0.0130094690177091
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 5.900730666666667
$ws.Cells.Item(2,14).Value = 17.702192
$ws.Cells.Item(2,15).Value = 0.03970749001357476
$ws.Cells.Item(2,16).Value = 0.03970749001357476
$ws.Cells.Item(2,17).Value = 1.447124692346667
$ws.Cells.Item(2,18).Value = 13.02412223112
$ws.Cells.Item(2,19).Value = 0.0005165733611025943
$ws.Cells.Item(2,20).Value = 0.0005165733611025942

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Fgf7"
$ws.Cells.Item(3,3).Value = "Fgfr1"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 2
$ws.Cells.Item(3,6).Value = 0.6666666666666666
$ws.Cells.Item(3,7).Value = 0.245245
$ws.Cells.Item(3,8).Value = 0.735735
$ws.Cells.Item(3,9).Value = 0.0130094690177091
$ws.Cells.Item(3,10).Value = 0.0130094690177091
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 111.5917106666667
$ws.Cells.Item(3,14).Value = 334.775132
$ws.Cells.Item(3,15).Value = 0.7509284844884279
$ws.Cells.Item(3,16).Value = 0.7509284844884279
$ws.Cells.Item(3,17).Value = 27.36730908244667
$ws.Cells.Item(3,18).Value = 246.30578174202
$ws.Cells.Item(3,19).Value = 0.00976918085346745
$ws.Cells.Item(3,20).Value = 0.00976918085346745

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Fgf7"
$ws.Cells.Item(4,3).Value = "Fgfr1"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 2
$ws.Cells.Item(4,6).Value = 0.6666666666666666
$ws.Cells.Item(4,7).Value = 0.245245
$ws.Cells.Item(4,8).Value = 0.735735
$ws.Cells.Item(4,9).Value = 0.0130094690177091
$ws.Cells.Item(4,10).Value = 0.0130094690177091
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 31.11253633333333
$ws.Cells.Item(4,14).Value = 93.337609
$ws.Cells.Item(4,15).Value = 0.2093640254979974
$ws.Cells.Item(4,16).Value = 0.2093640254979974
$ws.Cells.Item(4,17).Value = 7.630193973068335
$ws.Cells.Item(4,18).Value = 68.671745757615
$ws.Cells.Item(4,19).Value = 0.002723714803139056
$ws.Cells.Item(4,20).Value = 0.002723714803139055

# Row 5
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Fgf7"
$ws.Cells.Item(5,3).Value = "Fgfr1"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 17.39906333333333
$ws.Cells.Item(5,8).Value = 52.19719000000001
$ws.Cells.Item(5,9).Value = 0.9229650976458578
$ws.Cells.Item(5,10).Value = 0.9229650976458579
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 5.900730666666667
$ws.Cells.Item(5,14).Value = 17.702192
$ws.Cells.Item(5,15).Value = 0.03970749001357476
$ws.Cells.Item(5,16).Value = 0.03970749001357476
$ws.Cells.Item(5,17).Value = 102.6671865822756
$ws.Cells.Item(5,18).Value = 924.0046792404801
$ws.Cells.Item(5,19).Value = 0.03664862739765096
$ws.Cells.Item(5,20).Value = 0.03664862739765095

# Row 6
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Fgf7"
$ws.Cells.Item(6,3).Value = "Fgfr1"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 17.39906333333333
$ws.Cells.Item(6,8).Value = 52.19719000000001
$ws.Cells.Item(6,9).Value = 0.9229650976458578
$ws.Cells.Item(6,10).Value = 0.9229650976458579
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 111.5917106666667
$ws.Cells.Item(6,14).Value = 334.775132
$ws.Cells.Item(6,15).Value = 0.7509284844884279
$ws.Cells.Item(6,16).Value = 0.7509284844884279
$ws.Cells.Item(6,17).Value = 1941.591241364342
$ws.Cells.Item(6,18).Value = 17474.32117227908
$ws.Cells.Item(6,19).Value = 0.6930807820109178
$ws.Cells.Item(6,20).Value = 0.6930807820109179

# Row 7
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Fgf7"
$ws.Cells.Item(7,3).Value = "Fgfr1"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 17.39906333333333
$ws.Cells.Item(7,8).Value = 52.19719000000001
$ws.Cells.Item(7,9).Value = 0.9229650976458578
$ws.Cells.Item(7,10).Value = 0.9229650976458579
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 31.11253633333333
$ws.Cells.Item(7,14).Value = 93.337609
$ws.Cells.Item(7,15).Value = 0.2093640254979974
$ws.Cells.Item(7,16).Value = 0.2093640254979974
$ws.Cells.Item(7,17).Value = 541.3289901243012
$ws.Cells.Item(7,18).Value = 4871.960911118711
$ws.Cells.Item(7,19).Value = 0.1932356882372891
$ws.Cells.Item(7,20).Value = 0.193235688237289

# Row 8
$ws.Cells.Item(8,1).Value = "sCs"
$ws.Cells.Item(8,2).Value = "Fgf7"
$ws.Cells.Item(8,3).Value = "Fgfr1"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 1.206960666666667
$ws.Cells.Item(8,8).Value = 3.620882
$ws.Cells.Item(8,9).Value = 0.06402543333643303
$ws.Cells.Item(8,10).Value = 0.06402543333643303
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 5.900730666666667
$ws.Cells.Item(8,14).Value = 17.702192
$ws.Cells.Item(8,15).Value = 0.03970749001357476
$ws.Cells.Item(8,16).Value = 0.03970749001357476
$ws.Cells.Item(8,17).Value = 7.121949819260444
$ws.Cells.Item(8,18).Value = 64.097548373344
$ws.Cells.Item(8,19).Value = 0.002542289254821211
$ws.Cells.Item(8,20).Value = 0.002542289254821211

# Row 9
$ws.Cells.Item(9,1).Value = "sCs"
$ws.Cells.Item(9,2).Value = "Fgf7"
$ws.Cells.Item(9,3).Value = "Fgfr1"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 1.206960666666667
$ws.Cells.Item(9,8).Value = 3.620882
$ws.Cells.Item(9,9).Value = 0.06402543333643303
$ws.Cells.Item(9,10).Value = 0.06402543333643303
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 111.5917106666667
$ws.Cells.Item(9,14).Value = 334.775132
$ws.Cells.Item(9,15).Value = 0.7509284844884279
$ws.Cells.Item(9,16).Value = 0.7509284844884279
$ws.Cells.Item(9,17).Value = 134.6868055007137
$ws.Cells.Item(9,18).Value = 1212.181249506424
$ws.Cells.Item(9,19).Value = 0.04807852162404252
$ws.Cells.Item(9,20).Value = 0.04807852162404252

# Row 10
$ws.Cells.Item(10,1).Value = "sCs"
$ws.Cells.Item(10,2).Value = "Fgf7"
$ws.Cells.Item(10,3).Value = "Fgfr1"
$ws.Cells.Item(10,4).Value = "sCs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 1.206960666666667
$ws.Cells.Item(10,8).Value = 3.620882
$ws.Cells.Item(10,9).Value = 0.06402543333643303
$ws.Cells.Item(10,10).Value = 0.06402543333643303
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 31.11253633333333
$ws.Cells.Item(10,14).Value = 93.337609
$ws.Cells.Item(10,15).Value = 0.2093640254979974
$ws.Cells.Item(10,16).Value = 0.2093640254979974
$ws.Cells.Item(10,17).Value = 37.55160759457089
$ws.Cells.Item(10,18).Value = 337.964468351138
$ws.Cells.Item(10,19).Value = 0.0134046224575693
$ws.Cells.Item(10,20).Value = 0.0134046224575693

